$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (shifts existing A:D to B:E), preserving their widths/content
$ws.Columns.Item(1).Insert()

# Row 1 (headers)
$ws.Cells.Item(1,1).Value = 'TabName'
$ws.Cells.Item(1,2).Value = 'query'
$ws.Cells.Item(1,3).Value = 'StatQuery'
$ws.Cells.Item(1,4).Value = 'dbExcel'
$ws.Cells.Item(1,5).Value = 'WebExcel'

# Row 2 (data)
$ws.Cells.Item(2,1).Value = 'CasesTab'
$ws.Cells.Item(2,2).Value = 'MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = "AMERICAN_INDIAN_OR_ALASKA_NATIVE"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '''') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '''') AS `Trial Code`,
    COALESCE(a.arm_id, '''') AS `Arm`,
    COALESCE(a.arm_drug, '''') AS `Arm Treatment`,
    COALESCE(c.disease, '''') AS `Diagnosis`,
    COALESCE(c.gender, '''') AS `Gender`,
    COALESCE(c.race, '''') AS `Race`,
    COALESCE(c.ethnicity, '''') AS `Ethnicity`'
$ws.Cells.Item(2,3).Value = 'MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = "AMERICAN_INDIAN_OR_ALASKA_NATIVE"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials'
$ws.Cells.Item(2,4).Value = 'TC01_Trials_Filter_Race-AmerIndAlask_Neo4jData.xlsx'
$ws.Cells.Item(2,5).Value = 'TC01_Trials_Filter_Race-AmerIndAlask_WebData.xlsx'

# Wrap text for the long query cells (matches style index 1 used previously)
$ws.Cells.Item(2,2).WrapText = $true
$ws.Cells.Item(2,3).WrapText = $true

# Column A width (narrow, best-fit style column for short tab names)
$ws.Columns.Item(1).ColumnWidth = 8

# Row height to fit the new multi-line query text
$ws.Rows.Item(2).RowHeight = 174

# Selection as recorded after the edit
$ws.Range("B4").Select()
